$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.756.24"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").Value = "1.604.25"
$ws.Range("E3").Value = "  +0.49%  "
$ws.Range("D4").Value = "1.01"
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").Value = "212.03"
$ws.Range("E5").Value = "  +0.28%  "
$ws.Range("D6").Value = "0.514"
$ws.Range("E6").Value = "  +0.22%  "
$ws.Range("D7").Value = "1.01"
$ws.Range("E7").Value = "  +0.24%  "
$ws.Range("D8").Value = "0.0620"
$ws.Range("E8").Value = "  +0.26%  "
$ws.Range("D9").Value = "0.247"
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("D10").Value = "19.64"
$ws.Range("E10").Value = "  +0.70%  "
$ws.Range("D11").Value = "0.0848"
$ws.Range("E11").Value = "  +0.82%  "
$ws.Range("D12").Value = "1.828.68"
$ws.Range("E12").Value = "  +0.43%  "
$ws.Range("D13").Value = "1.592.12"
$ws.Range("E13").Value = "  +0.19%  "
$ws.Range("D14").Value = "4.07"
$ws.Range("E14").Value = "  +1.10%  "
$ws.Range("D15").Value = "0.526"
$ws.Range("E15").Value = "  +0.52%  "
$ws.Range("D16").Value = "65.09"
$ws.Range("E16").Value = "  +0.12%  "
$ws.Range("D17").Value = "0.0₃0742"
$ws.Range("E17").Value = "  -0.22%  "
$ws.Range("D18").Value = "209.38"
$ws.Range("E18").Value = "  +0.12%  "
$ws.Range("D19").Value = "1.01"
$ws.Range("E19").Value = "  +0.26%  "
$ws.Range("D20").Value = "7.16"
$ws.Range("E20").Value = "  +1.42%  "
$ws.Range("D21").Value = "4.31"
$ws.Range("E21").Value = "  +0.34%  "
$ws.Range("D22").Value = "2.22"
$ws.Range("E22").Value = "  -4.71%  "
$ws.Range("D23").Value = "9.04"
$ws.Range("E23").Value = "  +0.72%  "
$ws.Range("D24").Value = "143.84"
$ws.Range("E24").Value = "  +0.52%  "
$ws.Range("D25").Value = "1.01"
$ws.Range("E25").Value = "  +0.40%  "
$ws.Range("D26").Value = "7.09"
$ws.Range("E26").Value = "  -0.27%  "
$ws.Range("D27").Value = "0.114"
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").Value = "15.37"
$ws.Range("E28").Value = "  +0.31%  "
$ws.Range("D29").Value = "0.0507"
$ws.Range("E29").Value = "  -1.33%  "
$ws.Range("D30").Value = "1.16"
$ws.Range("E30").Value = "  +0.44%  "
$ws.Range("D31").Value = "3.28"
$ws.Range("E31").Value = "  +1.17%  "
$ws.Range("D32").Value = "2.96"
$ws.Range("E32").Value = "  +0.60%  "
$ws.Range("D33").Value = "1.289.67"
$ws.Range("E33").Value = "  +0.10%  "
$ws.Range("D34").Value = "2.48"
$ws.Range("E34").Value = "  +1.24%  "
$ws.Range("D35").Value = "1.22"
$ws.Range("E35").Value = "  +16.40%  "
$ws.Range("D36").Value = "1.49"
$ws.Range("E36").Value = "  +0.48%  "
$ws.Range("D37").Value = "0.587"
$ws.Range("E37").Value = "  -5.16%  "
$ws.Range("D38").Value = "0.0170"
$ws.Range("E38").Value = "  -0.38%  "
$ws.Range("D39").Value = "0.826"
$ws.Range("E39").Value = "  -0.29%  "
$ws.Range("D40").Value = "5.45"
$ws.Range("E40").Value = "  -0.13%  "
$ws.Range("D41").Value = "2.19"
$ws.Range("E41").Value = "  +0.07%  "
$ws.Range("D42").Value = "0.779"
$ws.Range("E42").Value = "  -0.28%  "
$ws.Range("D43").Value = "62.54"
$ws.Range("E43").Value = "  -1.09%  "
$ws.Range("D44").Value = "1.740.38"
$ws.Range("E44").Value = "  +0.50%  "
$ws.Range("D45").Value = "90.52"
$ws.Range("E45").Value = "  -0.66%  "
$ws.Range("D46").Value = "1.57"
$ws.Range("E46").Value = "  +0.33%  "
$ws.Range("D47").Value = "0.102"
$ws.Range("E47").Value = "  +1.04%  "
$ws.Range("D48").Value = "0.0513"
$ws.Range("E48").Value = "  +0.72%  "
$ws.Range("D49").Value = "7.57"
$ws.Range("E49").Value = "  +2.95%  "
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  +0.15%  "
$ws.Range("D51").Value = "0.400"
$ws.Range("E51").Value = "  +1.91%  "
